$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('2025-12-08 04:26:12', 'Admin', 'Login', 'login_success', 'Role: admin'),
    @('2025-12-08 04:26:12', 'Admin', 'dashboard', 'access_granted', 'Opened dashboard page'),
    @('2025-12-08 04:26:14', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 04:34:16', 'Admin', 'Login', 'login_success', 'Role: admin'),
    @('2025-12-08 04:34:16', 'Admin', 'dashboard', 'access_granted', 'Opened dashboard page'),
    @('2025-12-08 04:34:18', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 04:34:22', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 04:34:25', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 04:36:13', 'Admin', 'Login', 'login_success', 'Role: admin'),
    @('2025-12-08 04:36:13', 'Admin', 'dashboard', 'access_granted', 'Opened dashboard page'),
    @('2025-12-08 04:36:15', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 04:44:04', 'Admin', 'Login', 'login_success', 'Role: admin'),
    @('2025-12-08 04:44:04', 'Admin', 'dashboard', 'access_granted', 'Opened dashboard page'),
    @('2025-12-08 04:45:24', 'Admin', 'Login', 'login_success', 'Role: admin'),
    @('2025-12-08 04:45:24', 'Admin', 'dashboard', 'access_granted', 'Opened dashboard page'),
    @('2025-12-08 04:48:07', 'Admin', 'receipt', 'access_granted', 'Opened receipt page'),
    @('2025-12-08 04:48:08', 'Admin', 'customers', 'access_granted', 'Opened customers page'),
    @('2025-12-08 04:48:09', 'Admin', 'dashboard', 'access_granted', 'Opened dashboard page'),
    @('2025-12-08 04:48:10', 'Admin', 'dashboard', 'access_granted', 'Opened dashboard page'),
    @('2025-12-08 04:48:21', 'Admin', 'settings', 'access_granted', 'Opened settings page'),
    @('2025-12-08 04:55:13', 'Admin', 'Login', 'login_success', 'Role: admin'),
    @('2025-12-08 04:55:13', 'Admin', 'dashboard', 'access_granted', 'Opened dashboard page'),
    @('2025-12-08 04:55:15', 'Admin', 'dashboard', 'access_granted', 'Opened dashboard page'),
    @('2025-12-08 04:55:16', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 04:55:17', 'Admin', 'invoice', 'access_granted', 'Opened invoice page'),
    @('2025-12-08 04:55:46', 'Admin', 'dashboard', 'access_granted', 'Opened dashboard page'),
    @('2025-12-08 04:55:50', 'Admin', 'receipt', 'access_granted', 'Opened receipt page'),
    @('2025-12-08 04:55:51', 'Admin', 'customers', 'access_granted', 'Opened customers page'),
    @('2025-12-08 04:55:53', 'Admin', 'settings', 'access_granted', 'Opened settings page'),
    @('2025-12-08 04:56:08', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 04:56:13', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 04:56:19', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 04:57:02', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 04:57:09', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 04:57:15', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 04:57:22', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 04:57:29', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 05:15:42', 'Admin', 'Login', 'login_success', 'Role: admin'),
    @('2025-12-08 05:15:42', 'Admin', 'dashboard', 'access_granted', 'Opened dashboard page'),
    @('2025-12-08 05:15:45', 'Admin', 'dashboard', 'access_granted', 'Opened dashboard page'),
    @('2025-12-08 05:15:46', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 05:15:50', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 05:15:52', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 05:20:59', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 05:21:03', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 05:22:12', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 05:22:13', 'Admin', 'dashboard', 'access_granted', 'Opened dashboard page'),
    @('2025-12-08 05:28:44', 'Admin', 'Login', 'login_success', 'Role: admin'),
    @('2025-12-08 05:28:44', 'Admin', 'dashboard', 'access_granted', 'Opened dashboard page'),
    @('2025-12-08 05:28:46', 'Admin', 'dashboard', 'access_granted', 'Opened dashboard page'),
    @('2025-12-08 05:28:47', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 05:28:51', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 05:28:52', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 05:30:39', 'Admin', 'Login', 'login_success', 'Role: admin'),
    @('2025-12-08 05:30:39', 'Admin', 'dashboard', 'access_granted', 'Opened dashboard page'),
    @('2025-12-08 05:30:49', 'Admin', 'settings', 'access_granted', 'Opened settings page'),
    @('2025-12-08 05:30:53', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 05:30:55', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 05:30:55', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 05:30:56', 'Admin', 'quotation', 'access_granted', 'Opened quotation page'),
    @('2025-12-08 05:30:59', 'Admin', 'quotation', 'access_granted', 'Opened quotation page')
)

$startRow = 226
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}

Write-Host "Added $($data.Count) rows starting at row $startRow"
